$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1150.5217
$ws.Range("I28").Value = 1229.7646
$ws.Range("K28").Value = 1229.7646
$ws.Range("M28").Value = -744.7646

$ws.Range("H62").Value = 2535.7896
$ws.Range("I62").Value = 2465.077
$ws.Range("K62").Value = 2465.077
$ws.Range("M62").Value = -1841.077

$ws.Range("H65").Value = 2535.7896
$ws.Range("I65").Value = 2465.077
$ws.Range("K65").Value = 12325.385
$ws.Range("M65").Value = -9205.385000000002

$ws.Range("H132").Value = 3579.1738
$ws.Range("I132").Value = 3579.1738
$ws.Range("K132").Value = 10737.5214
$ws.Range("M132").Value = -8207.5214

$ws.Range("H137").Value = 1462.8
$ws.Range("I137").Value = 1279.0526
$ws.Range("J137").Value = 1629.0476
$ws.Range("K137").Value = 3837.1578
$ws.Range("L137").Value = 4887.142800000001
$ws.Range("M137").Value = -1287.1578
$ws.Range("N137").Value = -9987.142800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1585.9333
$ws.Range("I2").Value = 1606.76
$ws.Range("J2").Value = 1481.8
$ws.Range("K2").Value = 1606.76
$ws.Range("L2").Value = 1481.8
$ws.Range("M2").Value = -1493.76
$ws.Range("N2").Value = -1707.8

$ws.Range("H61").Value = 3570.3333
$ws.Range("J61").Value = 4440.6665
$ws.Range("L61").Value = 4440.6665
$ws.Range("N61").Value = -4864.6665

$ws.Range("H74").Value = 17594.4
$ws.Range("I74").Value = 1030.641
$ws.Range("K74").Value = 1030.641
$ws.Range("M74").Value = -156.6410000000001

$ws.Range("H77").Value = 17594.4
$ws.Range("I77").Value = 1030.641
$ws.Range("K77").Value = 5153.205
$ws.Range("M77").Value = -785.2049999999999

$ws.Range("H116").Value = 1585.9333
$ws.Range("I116").Value = 1606.76
$ws.Range("J116").Value = 1481.8
$ws.Range("K116").Value = 1606.76
$ws.Range("L116").Value = 1481.8
$ws.Range("M116").Value = 687.24
$ws.Range("N116").Value = -6069.8

$ws.Range("H132").Value = 2923.578
$ws.Range("I132").Value = 2766.923
$ws.Range("K132").Value = 8300.769
$ws.Range("M132").Value = -5770.769

$ws.Range("H136").Value = 3570.3333
$ws.Range("J136").Value = 4440.6665
$ws.Range("L136").Value = 13321.9995
$ws.Range("N136").Value = -18421.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1585.9333
$ws.Range("I3").Value = 1606.76
$ws.Range("J3").Value = 1481.8
$ws.Range("K3").Value = 1606.76
$ws.Range("L3").Value = 1481.8
$ws.Range("M3").Value = -1492.76
$ws.Range("N3").Value = -1709.8

$ws.Range("H107").Value = 5301.5
$ws.Range("I107").Value = 5482.6665
$ws.Range("J107").Value = 1497
$ws.Range("K107").Value = 5482.6665
$ws.Range("L107").Value = 1497
$ws.Range("M107").Value = -3562.6665
$ws.Range("N107").Value = -5337

$ws.Range("H134").Value = 2311.804
$ws.Range("I134").Value = 2283.449
$ws.Range("J134").Value = 3006.5
$ws.Range("K134").Value = 6850.347
$ws.Range("L134").Value = 9019.5
$ws.Range("M134").Value = -4315.347
$ws.Range("N134").Value = -14089.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 5329.5
$ws.Range("I14").Value = 209.5
$ws.Range("K14").Value = 209.5
$ws.Range("M14").Value = -39.5

$ws.Range("H22").Value = 1151.0667
$ws.Range("I22").Value = 634.8
$ws.Range("K22").Value = 634.8
$ws.Range("M22").Value = -284.8

$ws.Range("H25").Value = 1732.6154
$ws.Range("I25").Value = 1002.2727
$ws.Range("K25").Value = 1002.2727
$ws.Range("M25").Value = -828.2727

$ws.Range("H31").Value = 2303.3088
$ws.Range("I31").Value = 1703.5264
$ws.Range("K31").Value = 1703.5264
$ws.Range("M31").Value = -1408.5264

$ws.Range("H34").Value = 2303.3088
$ws.Range("I34").Value = 1703.5264
$ws.Range("K34").Value = 1703.5264
$ws.Range("M34").Value = -1501.5264

$ws.Range("H99").Value = 3338.0715
$ws.Range("I99").Value = 3342.037
$ws.Range("K99").Value = 3342.037
$ws.Range("M99").Value = -1844.037

$ws.Range("H126").Value = 3338.0715
$ws.Range("I126").Value = 3342.037
$ws.Range("K126").Value = 10026.111
$ws.Range("M126").Value = -7556.110999999999

$ws.Range("H132").Value = 8981
$ws.Range("I132").Value = 7962.3335
$ws.Range("K132").Value = 23887.0005
$ws.Range("M132").Value = -21357.0005

$ws.Range("H134").Value = 3715.96
$ws.Range("I134").Value = 3917
$ws.Range("J134").Value = 2660.5
$ws.Range("K134").Value = 11751
$ws.Range("L134").Value = 7981.5
$ws.Range("M134").Value = -9216
$ws.Range("N134").Value = -13051.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 5500
$ws.Range("J17").Value = 5500
$ws.Range("L17").Value = 16500
$ws.Range("N17").Value = -16838

$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 8996.666999999999
$ws.Range("J21").Value = 9245
$ws.Range("L21").Value = 9245
$ws.Range("N21").Value = -9591

$ws.Range("H30").Value = 8996.666999999999
$ws.Range("J30").Value = 9245
$ws.Range("L30").Value = 9245
$ws.Range("N30").Value = -9455

$ws.Range("H102").Value = 26984.783
$ws.Range("I102").Value = 42600.6
$ws.Range("J102").Value = 8394.522999999999
$ws.Range("K102").Value = 42600.6
$ws.Range("L102").Value = 8394.522999999999
$ws.Range("M102").Value = -40978.6
$ws.Range("N102").Value = -11638.523

$ws.Range("H122").Value = 83399.89
$ws.Range("I122").Value = 109706.29
$ws.Range("J122").Value = 1557.7778
$ws.Range("K122").Value = 329118.87
$ws.Range("L122").Value = 4673.3334
$ws.Range("M122").Value = -326668.87
$ws.Range("N122").Value = -9573.3334

$ws.Range("H126").Value = 115513.75
$ws.Range("I126").Value = 131301.42
$ws.Range("K126").Value = 393904.26
$ws.Range("M126").Value = -391434.26

$ws.Range("H132").Value = 2864.5925
$ws.Range("I132").Value = 2905.84
$ws.Range("J132").Value = 2349
$ws.Range("K132").Value = 8717.52
$ws.Range("L132").Value = 7047
$ws.Range("M132").Value = -6187.52
$ws.Range("N132").Value = -12107

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12358.818
$ws.Range("I7").Value = 15094.875
$ws.Range("J7").Value = 5062.6665
$ws.Range("K7").Value = 15094.875
$ws.Range("L7").Value = 5062.6665
$ws.Range("M7").Value = -14982.875
$ws.Range("N7").Value = -5286.6665

$ws.Range("H40").Value = 3768.5518
$ws.Range("I40").Value = 3687.84
$ws.Range("J40").Value = 4273
$ws.Range("K40").Value = 3687.84
$ws.Range("L40").Value = 4273
$ws.Range("M40").Value = -3551.84
$ws.Range("N40").Value = -4545

$ws.Range("H46").Value = 1848.7
$ws.Range("I46").Value = 678.3333
$ws.Range("K46").Value = 678.3333
$ws.Range("M46").Value = -490.3333

$ws.Range("H93").Value = 2278.7112
$ws.Range("I93").Value = 1944.9062
$ws.Range("J93").Value = 3100.3845
$ws.Range("K93").Value = 1944.9062
$ws.Range("L93").Value = 3100.3845
$ws.Range("M93").Value = -696.9061999999999
$ws.Range("N93").Value = -5596.3845

$ws.Range("H126").Value = 12358.818
$ws.Range("I126").Value = 15094.875
$ws.Range("J126").Value = 5062.6665
$ws.Range("K126").Value = 45284.625
$ws.Range("L126").Value = 15187.9995
$ws.Range("M126").Value = -42814.625
$ws.Range("N126").Value = -20127.9995

$ws.Range("H136").Value = 40211.668
$ws.Range("I136").Value = 2477.9473
$ws.Range("K136").Value = 7433.841899999999
$ws.Range("M136").Value = -4883.841899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 43481610
$ws.Range("I126").Value = 71432340
$ws.Range("K126").Value = 214297020
$ws.Range("M126").Value = -214294550

$ws.Range("H132").Value = 3473.4443
$ws.Range("I132").Value = 3570.3125
$ws.Range("J132").Value = 2698.5
$ws.Range("K132").Value = 10710.9375
$ws.Range("L132").Value = 8095.5
$ws.Range("M132").Value = -8180.9375
$ws.Range("N132").Value = -13155.5

$ws.Range("H136").Value = 51285224
$ws.Range("J136").Value = 76926110
$ws.Range("L136").Value = 230778330
$ws.Range("N136").Value = -230783430
